# Refresh the crypto price/volume table with the latest scraped values.
# Price/Volume cells are stored as literal text (not numbers) in the sheet,
# so for values that look like a plain decimal number we force the cell's
# NumberFormat to Text ("@") first - otherwise Excel would silently convert
# a value like "519.97" into a floating point number (losing exact text,
# e.g. "1.00" -> 1) instead of keeping it as the original text string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.171.25'
$ws.Range('E2').Value = '  -1.78%  '
$ws.Range('D3').Value = '2.471.12'
$ws.Range('E3').Value = '  -2.15%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '519.97'
$ws.Range('E5').Value = '  -3.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.67'
$ws.Range('E6').Value = '  -3.90%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.558'
$ws.Range('E8').Value = '  -1.71%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0994'
$ws.Range('E9').Value = '  -2.24%  '
$ws.Range('E10').Value = '  -0.61%  '
$ws.Range('E11').Value = '  +0.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.343'
$ws.Range('E12').Value = '  -2.02%  '
$ws.Range('D13').Value = '2.907.41'
$ws.Range('E13').Value = '  -2.25%  '
$ws.Range('D14').Value = '58.099.37'
$ws.Range('E14').Value = '  -1.80%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '22.15'
$ws.Range('E15').Value = '  -4.39%  '
$ws.Range('E16').Value = '  -2.46%  '
$ws.Range('D17').Value = '2.469.06'
$ws.Range('E17').Value = '  -2.29%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.87'
$ws.Range('E18').Value = '  -2.34%  '
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.19'
$ws.Range('E19').Value = '  -2.60%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '320.77'
$ws.Range('E20').Value = '  -1.58%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.00'
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.76'
$ws.Range('E22').Value = '  -3.87%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '64.42'
$ws.Range('E23').Value = '  -2.08%  '
$ws.Range('E24').Value = '  -3.63%  '
$ws.Range('E25').Value = '  -0.26%  '
$ws.Range('E26').Value = '  -3.55%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.42'
$ws.Range('E27').Value = '  -3.47%  '
$ws.Range('D28').Value = '0.0₃0751'
$ws.Range('E28').Value = '  -3.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.41'
$ws.Range('E29').Value = '  -5.14%  '
$ws.Range('E30').Value = '  -4.77%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '165.70'
$ws.Range('E31').Value = '  +1.56%  '
$ws.Range('E32').Value = '  -4.13%  '
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('E34').Value = '  -0.19%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '18.15'
$ws.Range('E35').Value = '  -1.88%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.34'
$ws.Range('E36').Value = '  -9.32%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.00'
$ws.Range('E37').Value = '  -3.17%  '
$ws.Range('E38').Value = '  -3.96%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.797'
$ws.Range('E39').Value = '  -2.82%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '276.99'
$ws.Range('E40').Value = '  -3.66%  '
$ws.Range('E41').Value = '  -4.92%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.03'
$ws.Range('E42').Value = '  -3.82%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.594'
$ws.Range('E43').Value = '  -2.81%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '126.54'
$ws.Range('E44').Value = '  -4.48%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0909'
$ws.Range('E45').Value = '  -2.58%  '
$ws.Range('E46').Value = '  -3.64%  '
$ws.Range('E47').Value = '  -3.52%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '17.18'
$ws.Range('E48').Value = '  -1.46%  '
$ws.Range('D49').Value = '1.736.61'
$ws.Range('E49').Value = '  -1.53%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.67'
$ws.Range('E51').Value = '  -2.03%  '
